# ------------------------------------------------------------------------
# Applies the "additional scraping" edit:
#  1. Adds a new "Player Info" sheet (before the existing one) with the
#     player's ID/NAME/BATTING_HAND/BOWL_STYLE.
#  2. Renames the MATCH_CARD_LINK column in the existing "ODI Batting"
#     sheet to MATCH_CODE and replaces the full scorecard URL with just
#     the numeric match code.
#  3. Adds a new "ODI Batting Extra" sheet (after "ODI Batting") with
#     additional per-match batting stats.
# ------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-HeaderCell($ws, $cell, $text) {
    $ws.Range($cell).Value = $text
    $ws.Range($cell).Font.Bold = $true
    $ws.Range($cell).HorizontalAlignment = -4108   # xlCenter
    $ws.Range($cell).VerticalAlignment = -4160     # xlTop
    $ws.Range($cell).Borders.LineStyle = 1
}

function Set-StandardMargins($ws) {
    # Matches the pageMargins already used by the original "ODI Batting"
    # sheet (0.75in / 1in / 0.5in, expressed in points).
    $ws.PageSetup.LeftMargin = 54
    $ws.PageSetup.RightMargin = 54
    $ws.PageSetup.TopMargin = 72
    $ws.PageSetup.BottomMargin = 72
    $ws.PageSetup.HeaderMargin = 36
    $ws.PageSetup.FooterMargin = 36
}

# ------------------------------------------------------------------------
# 1. "Player Info" sheet - inserted before the existing "ODI Batting" sheet
# ------------------------------------------------------------------------
$wsOdi = $wb.Worksheets.Item("ODI Batting")
$wsInfo = $wb.Worksheets.Add($wsOdi)
$wsInfo.Name = "Player Info"
Set-StandardMargins $wsInfo

$wsInfo.Range("A1:D2").NumberFormat = "@"

Set-HeaderCell $wsInfo "A1" "ID"
Set-HeaderCell $wsInfo "B1" "NAME"
Set-HeaderCell $wsInfo "C1" "BATTING_HAND"
Set-HeaderCell $wsInfo "D1" "BOWL_STYLE"

$wsInfo.Range("A2").Value = "4726"
$wsInfo.Range("B2").Value = "Alex Tyson Carey"
$wsInfo.Range("C2").Value = "Left Handed"
$wsInfo.Range("D2").Value = "Left Arm Medium Fast"

# ------------------------------------------------------------------------
# 2. "ODI Batting" sheet - rename MATCH_CARD_LINK -> MATCH_CODE and
#    collapse the scorecard URL down to the bare match code.
# ------------------------------------------------------------------------
# NOTE: inserting a new worksheet (above) invalidates previously captured
# worksheet handles (they keep pointing at the original *index*, not the
# original *sheet*), so re-resolve "ODI Batting" by name before using it.
$wsOdi = $wb.Worksheets.Item("ODI Batting")

$wsOdi.Range("D1").Value = "MATCH_CODE"

$lastRow = $wsOdi.Cells.Item($wsOdi.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 65 }

# Keep column D as text (it already held the text URL, and the replacement
# match-code values must not silently become numbers).
$wsOdi.Range("D2:D$lastRow").NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsOdi.Cells.Item($r, 4)
    $link = $cell.Value()
    if ($link -ne $null -and $link -ne "") {
        $idx = $link.LastIndexOf("MatchCode=")
        if ($idx -ge 0) {
            $code = $link.Substring($idx + 10)
            $cell.Value = $code
        }
    }
}

# ------------------------------------------------------------------------
# 3. "ODI Batting Extra" sheet - inserted after "ODI Batting"
# ------------------------------------------------------------------------
# Re-fetch the "ODI Batting" worksheet reference again (same staleness
# caveat as above) right before using it as the "After" anchor.
$wsOdi = $wb.Worksheets.Item("ODI Batting")
$wsExtra = $wb.Worksheets.Add([System.Type]::Missing, $wsOdi)
$wsExtra.Name = "ODI Batting Extra"
Set-StandardMargins $wsExtra

$extraRows = @(
    @("4486", $null, $null, $null, $null, "NO"),
    @("4564", 6, "0", "0", "1.28%", "NO"),
    @("4565", $null, $null, $null, $null, "NO"),
    @("4567", 6, "6", "1", "26.67%", "NO"),
    @("4594", $null, $null, $null, $null, "NO"),
    @("4597", 6, "1", "0", "7.94%", "NO"),
    @("4600", $null, $null, $null, $null, "NO"),
    @("4601", 5, "2", "0", "7.48%", "NO"),
    @("4603", 6, "1", "0", "27.44%", "NO"),
    @("4644", 4, "2", "0", "4.98%", "NO"),
    @("4645", 4, "2", "0", "26.00%", "NO"),
    @("4646", 4, "1", "0", "2.84%", "NO"),
    @("4647", $null, $null, $null, $null, "NO"),
    @("4648", 6, "0", "0", "6.15%", "NO"),
    @("4649", 5, "3", "0", "15.73%", "NO"),
    @("4660", $null, $null, $null, $null, "NO"),
    @("4663", $null, $null, $null, $null, "NO"),
    @("4666", $null, $null, $null, $null, "NO"),
    @("4728", 6, $null, $null, $null, "NO"),
    @("4732", 6, "2", "1", "14.13%", "NO")
)

$extraLastRow = 1 + $extraRows.Count
$wsExtra.Range("A1:A$extraLastRow").NumberFormat = "@"
$wsExtra.Range("C1:F$extraLastRow").NumberFormat = "@"

Set-HeaderCell $wsExtra "A1" "MATCH_CODE"
Set-HeaderCell $wsExtra "B1" "BATTING_POSITION"
Set-HeaderCell $wsExtra "C1" "NUM_4"
Set-HeaderCell $wsExtra "D1" "NUM_6"
Set-HeaderCell $wsExtra "E1" "PERCENT_RUNS_OF_TOTAL"
Set-HeaderCell $wsExtra "F1" "MAN_OF_MATCH"

$r = 2
foreach ($row in $extraRows) {
    $wsExtra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) {
        $wsExtra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($row[2] -ne $null) {
        $wsExtra.Cells.Item($r, 3).Value = $row[2]
    }
    if ($row[3] -ne $null) {
        $wsExtra.Cells.Item($r, 4).Value = $row[3]
    }
    if ($row[4] -ne $null) {
        $wsExtra.Cells.Item($r, 5).Value = $row[4]
    }
    $wsExtra.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# ------------------------------------------------------------------------
# Make sure "ODI Batting" remains the active tab, matching the original
# workbook view (activeTab="0" referred to the sole original sheet).
# ------------------------------------------------------------------------
$wsOdi = $wb.Worksheets.Item("ODI Batting")
$wsOdi.Activate()
